$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell value updates per the crypto price/volume refresh (GitHub Actions bot).
# Cells whose new text is a bare number (e.g. "1.00", "7.60") need NumberFormat
# forced to text BEFORE the value is written, otherwise Excel auto-coerces the
# text into a numeric value (and mangles things like trailing zeros).

$ws.Range("D2").Value = '60.846.09'
$ws.Range("E2").Value = '  +6.50%  '
$ws.Range("D3").Value = '2.639.49'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '514.18'
$ws.Range("E5").Value = '  +5.10%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '159.56'
$ws.Range("E6").Value = '  +2.52%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.615'
$ws.Range("E7").Value = '  -0.92%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.996'
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("D9").Value = '2.688.97'
$ws.Range("E9").Value = '  +9.86%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.16'
$ws.Range("E10").Value = '  +8.35%  '
$ws.Range("E11").Value = '  +5.93%  '
$ws.Range("E12").Value = '  +4.15%  '
$ws.Range("E13").Value = '  +0.90%  '
$ws.Range("D14").Value = '3.125.89'
$ws.Range("E14").Value = '  +9.71%  '
$ws.Range("D15").Value = '61.056.30'
$ws.Range("E15").Value = '  +6.54%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '22.26'
$ws.Range("E16").Value = '  +6.80%  '
$ws.Range("E17").Value = '  +5.14%  '
$ws.Range("D18").Value = '2.683.93'
$ws.Range("E18").Value = '  +9.79%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.83'
$ws.Range("E19").Value = '  +1.11%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '348.92'
$ws.Range("E20").Value = '  +5.98%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.56'
$ws.Range("E21").Value = '  +5.72%  '
$ws.Range("E22").Value = '  +4.61%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.998'
$ws.Range("E23").Value = '  -0.11%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '60.64'
$ws.Range("E24").Value = '  +3.77%  '
$ws.Range("E25").Value = '  +3.63%  '
$ws.Range("D26").Value = '2.790.69'
$ws.Range("E26").Value = '  +10.00%  '
$ws.Range("E27").Value = '  +4.68%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.997'
$ws.Range("E28").Value = '  -0.02%  '
$ws.Range("D29").Value = '0.0₃0872'
$ws.Range("E29").Value = '  +10.33%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.60'
$ws.Range("E30").Value = '  +4.01%  '
$ws.Range("E31").Value = '  +0.17%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '19.69'
$ws.Range("E32").Value = '  +4.89%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '156.97'
$ws.Range("E33").Value = '  +5.28%  '
$ws.Range("E34").Value = '  +4.07%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.73'
$ws.Range("E35").Value = '  +7.60%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.10'
$ws.Range("E36").Value = '  +10.36%  '
$ws.Range("E37").Value = '  +6.44%  '
$ws.Range("E38").Value = '  +3.61%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.54'
$ws.Range("E39").Value = '  +11.82%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '310.01'
$ws.Range("E40").Value = '  +15.13%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.80'
$ws.Range("E41").Value = '  +7.59%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.838'
$ws.Range("E42").Value = '  +29.92%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '35.78'
$ws.Range("E43").Value = '  +4.38%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.649'
$ws.Range("E44").Value = '  +8.39%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0579'
$ws.Range("E45").Value = '  +8.20%  '
$ws.Range("E46").Value = '  -0.54%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '20.34'
$ws.Range("E47").Value = '  +15.67%  '
$ws.Range("B48").Value = 'RenderToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.04'
$ws.Range("E48").Value = '  +6.67%  '
$ws.Range("B49").Value = 'FirstDigitalUSD'
$ws.Range("C49").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.991'
$ws.Range("E49").Value = '  -0.63%  '
$ws.Range("E50").Value = '  +3.95%  '
$ws.Range("D51").Value = '2.041.51'
$ws.Range("E51").Value = '  +9.79%  '
